# Add two new expense rows below the existing one (row 1), turning the
# single-row sheet into a 3-row table (A1:D3). This backs the new "export
# to PDF" feature with a couple more sample rows to show off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns C (date) and D (amount) hold text that merely looks like a date /
# number ("2023-03-07", "1500.0"). Format those cells as Text first so
# Excel stores the literal strings instead of silently coercing them into
# a date serial / numeric value.
$ws.Range("C2:C3").NumberFormat = "@"
$ws.Range("D2:D3").NumberFormat = "@"

$ws.Range("A2").Value = "Food"
$ws.Range("B2").Value = "walkdms"
$ws.Range("C2").Value = "2023-03-07"
$ws.Range("D2").Value = "1500.0"

$ws.Range("A3").Value = "Transportation"
$ws.Range("B3").Value = "smws"
$ws.Range("C3").Value = "2023-03-07"
$ws.Range("D3").Value = "1234.0"
